# Generate Report for Handoff
#
# Rows 4-7 of the "zh-cn" and "de-de" sheets describe files that were
# "Ready for handoff" with Priority "low" and a stale "Latest Handoff
# Datetime". Re-generating the handoff report:
#   - promotes their Priority to "ht"
#   - refreshes "Latest Handoff Datetime" (column H) to the new
#     generation timestamp
# The "Overview" sheet's "Latest HO Xliff Generate Date" column mirrors
# the newest per-locale handoff timestamp (the de-de one, generated
# after zh-cn), so it is refreshed to match as well.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-27 22:31:54"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("H5").Value = "2016-08-27 22:31:54"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("H6").Value = "2016-08-27 22:31:54"
$zhcn.Range("E7").Value = "ht"
$zhcn.Range("H7").Value = "2016-08-27 22:31:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "ht"
$dede.Range("H4").Value = "2016-08-27 22:31:59"
$dede.Range("E5").Value = "ht"
$dede.Range("H5").Value = "2016-08-27 22:31:59"
$dede.Range("E6").Value = "ht"
$dede.Range("H6").Value = "2016-08-27 22:31:59"
$dede.Range("E7").Value = "ht"
$dede.Range("H7").Value = "2016-08-27 22:31:59"

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("G4").Value = "2016-08-27 22:31:59"
$ov.Range("G5").Value = "2016-08-27 22:31:59"
$ov.Range("G6").Value = "2016-08-27 22:31:59"
$ov.Range("G7").Value = "2016-08-27 22:31:59"
